$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1 "Divisão" reusing the same formatting (bold, centered, bordered)
# as the other header cells (B1/C1) by copying C1's formats onto D1.
$ws.Range("D1").Value = "Divisão"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill D2:D18 with "Primeira Divisão" for every row of data (plain, unstyled cells)
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 4).Value = "Primeira Divisão"
}
